$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 772.26746
$ws.Range("I15").Value = 772.26746
$ws.Range("K15").Value = 2316.80238
$ws.Range("M15").Value = -2147.80238
$ws.Range("H64").Value = 7999.5
$ws.Range("J64").Value = 7999.5
$ws.Range("L64").Value = 7999.5
$ws.Range("N64").Value = -8495.5
$ws.Range("H67").Value = 7999.5
$ws.Range("J67").Value = 7999.5
$ws.Range("L67").Value = 7999.5
$ws.Range("N67").Value = -9715.5
$ws.Range("H100").Value = 1240.7858
$ws.Range("I100").Value = 951.6923
$ws.Range("J100").Value = 4999
$ws.Range("K100").Value = 951.6923
$ws.Range("L100").Value = 4999
$ws.Range("M100").Value = -410.6923
$ws.Range("N100").Value = -6081
$ws.Range("H113").Value = 7115.9165
$ws.Range("J113").Value = 9798
$ws.Range("L113").Value = 9798
$ws.Range("N113").Value = -16306
$ws.Range("H137").Value = 64519644
$ws.Range("I137").Value = 41669556
$ws.Range("J137").Value = 142862800
$ws.Range("K137").Value = 125008668
$ws.Range("L137").Value = 428588400
$ws.Range("M137").Value = -125006118
$ws.Range("N137").Value = -428593500
$ws.Range("H138").Value = 5707.3687
$ws.Range("I138").Value = 2558.8572
$ws.Range("J138").Value = 6418.3228
$ws.Range("K138").Value = 7676.571599999999
$ws.Range("L138").Value = 19254.9684
$ws.Range("M138").Value = -2536.571599999999
$ws.Range("N138").Value = -29534.9684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38465156
$ws.Range("I61").Value = 47621028
$ws.Range("J61").Value = 10500
$ws.Range("K61").Value = 47621028
$ws.Range("L61").Value = 10500
$ws.Range("M61").Value = -47620816
$ws.Range("N61").Value = -10924
$ws.Range("H74").Value = 125142650
$ws.Range("I74").Value = 125142650
$ws.Range("K74").Value = 125142650
$ws.Range("M74").Value = -125141776
$ws.Range("H77").Value = 125142650
$ws.Range("I77").Value = 125142650
$ws.Range("K77").Value = 625713250
$ws.Range("M77").Value = -625708882
$ws.Range("H97").Value = 559.3077
$ws.Range("I97").Value = 355.33334
$ws.Range("J97").Value = 1416
$ws.Range("K97").Value = 355.33334
$ws.Range("L97").Value = 1416
$ws.Range("M97").Value = 140.66666
$ws.Range("N97").Value = -2408
$ws.Range("H110").Value = 19725.809
$ws.Range("I110").Value = 21044.578
$ws.Range("J110").Value = 7197.5
$ws.Range("K110").Value = 21044.578
$ws.Range("L110").Value = 7197.5
$ws.Range("M110").Value = -18999.578
$ws.Range("N110").Value = -11287.5
$ws.Range("H132").Value = 25649744
$ws.Range("I132").Value = 10159.621
$ws.Range("J132").Value = 100004536
$ws.Range("K132").Value = 30478.863
$ws.Range("L132").Value = 300013608
$ws.Range("M132").Value = -27948.863
$ws.Range("N132").Value = -300018668
$ws.Range("H136").Value = 38465156
$ws.Range("I136").Value = 47621028
$ws.Range("J136").Value = 10500
$ws.Range("K136").Value = 142863084
$ws.Range("L136").Value = 31500
$ws.Range("M136").Value = -142860534
$ws.Range("N136").Value = -36600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2468.158
$ws.Range("I20").Value = 4452
$ws.Range("J20").Value = 1310.9166
$ws.Range("K20").Value = 4452
$ws.Range("L20").Value = 1310.9166
$ws.Range("M20").Value = -4205
$ws.Range("N20").Value = -1804.9166
$ws.Range("H94").Value = 3217.476
$ws.Range("I94").Value = 1142.6111
$ws.Range("J94").Value = 15666.667
$ws.Range("K94").Value = 1142.6111
$ws.Range("L94").Value = 15666.667
$ws.Range("M94").Value = -691.6111000000001
$ws.Range("N94").Value = -16568.667
$ws.Range("H105").Value = 12125.111
$ws.Range("I105").Value = 15167.429
$ws.Range("J105").Value = 1477
$ws.Range("K105").Value = 15167.429
$ws.Range("L105").Value = 1477
$ws.Range("M105").Value = -13420.429
$ws.Range("N105").Value = -4971
$ws.Range("H134").Value = 2799.9524
$ws.Range("I134").Value = 3054.1177
$ws.Range("J134").Value = 1719.75
$ws.Range("K134").Value = 9162.3531
$ws.Range("L134").Value = 5159.25
$ws.Range("M134").Value = -6627.3531
$ws.Range("N134").Value = -10229.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4738.1
$ws.Range("I86").Value = 3473
$ws.Range("J86").Value = 5581.5
$ws.Range("K86").Value = 3473
$ws.Range("L86").Value = 5581.5
$ws.Range("M86").Value = -2350
$ws.Range("N86").Value = -7827.5
$ws.Range("H89").Value = 4738.1
$ws.Range("I89").Value = 3473
$ws.Range("J89").Value = 5581.5
$ws.Range("K89").Value = 17365
$ws.Range("L89").Value = 27907.5
$ws.Range("M89").Value = -11749
$ws.Range("N89").Value = -39139.5
$ws.Range("H99").Value = 9268.529
$ws.Range("I99").Value = 11243.25
$ws.Range("K99").Value = 11243.25
$ws.Range("M99").Value = -9745.25
$ws.Range("H105").Value = 8937.25
$ws.Range("I105").Value = 2845
$ws.Range("J105").Value = 35337
$ws.Range("K105").Value = 2845
$ws.Range("L105").Value = 35337
$ws.Range("M105").Value = -1098
$ws.Range("N105").Value = -38831
$ws.Range("H126").Value = 9268.529
$ws.Range("I126").Value = 11243.25
$ws.Range("K126").Value = 33729.75
$ws.Range("M126").Value = -31259.75
$ws.Range("H132").Value = 128992
$ws.Range("I132").Value = 225444.11
$ws.Range("J132").Value = 4982.143
$ws.Range("K132").Value = 676332.33
$ws.Range("L132").Value = 14946.429
$ws.Range("M132").Value = -673802.33
$ws.Range("N132").Value = -20006.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 424.65
$ws.Range("I22").Value = 417.4
$ws.Range("K22").Value = 1252.2
$ws.Range("M22").Value = -1083.2
$ws.Range("H27").Value = 424.65
$ws.Range("I27").Value = 417.4
$ws.Range("K27").Value = 1252.2
$ws.Range("M27").Value = -1150.2
$ws.Range("H105").Value = 3600
$ws.Range("J105").Value = 3600
$ws.Range("L105").Value = 10800
$ws.Range("N105").Value = -16042
$ws.Range("H139").Value = 3270.1667
$ws.Range("I139").Value = 2628.3572
$ws.Range("K139").Value = 7885.071599999999
$ws.Range("M139").Value = -2745.071599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 61726.656
$ws.Range("I70").Value = 78375.41
$ws.Range("K70").Value = 78375.41
$ws.Range("M70").Value = -78105.41
$ws.Range("H73").Value = 61726.656
$ws.Range("I73").Value = 78375.41
$ws.Range("K73").Value = 78375.41
$ws.Range("M73").Value = -77439.41
$ws.Range("H132").Value = 2014.2157
$ws.Range("I132").Value = 1916.575
$ws.Range("K132").Value = 5749.725
$ws.Range("M132").Value = -3219.725

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4840.5
$ws.Range("J61").Value = 6884.3335
$ws.Range("L61").Value = 6884.3335
$ws.Range("N61").Value = -7288.3335
$ws.Range("H93").Value = 1758.8125
$ws.Range("I93").Value = 1295.5
$ws.Range("K93").Value = 1295.5
$ws.Range("M93").Value = -47.5
$ws.Range("H113").Value = 4840.5
$ws.Range("J113").Value = 6884.3335
$ws.Range("L113").Value = 6884.3335
$ws.Range("N113").Value = -11224.3335
$ws.Range("H132").Value = 36367268
$ws.Range("I132").Value = 3410.575
$ws.Range("K132").Value = 10231.725
$ws.Range("M132").Value = -7701.724999999999
$ws.Range("H136").Value = 6691.2354
$ws.Range("I136").Value = 6161.7856
$ws.Range("K136").Value = 18485.3568
$ws.Range("M136").Value = -15935.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 21763526
$ws.Range("I122").Value = 33369500
$ws.Range("J122").Value = 2324.625
$ws.Range("K122").Value = 100108500
$ws.Range("L122").Value = 6973.875
$ws.Range("M122").Value = -100106050
$ws.Range("N122").Value = -11873.875
$ws.Range("H136").Value = 1469.5
$ws.Range("I136").Value = 931.1667
$ws.Range("K136").Value = 2793.5001
$ws.Range("M136").Value = -243.5001000000002
